$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date (serial 45406 -> 2024-04-24). Update it to 45432 (2024-05-20).
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# D29: 300 -> 955
$ws.Range("D29").Value = 955

# D30: 223.526 -> 376.8
$ws.Range("D30").Value = 376.8
